# outputs-HGR-r202/test-g__Amedibacillus_split_pruned.xlsx
# "updated outputs-r202, previous copy of ful-path.csv"
#
# The sheet holds a single quadratic-svm-score prediction row. The refreshed
# run re-scored the genome and produced a new (negative) decision-function
# value for B2; the header/label cells (A1:C1, A2) get their format
# reapplied as part of the re-export, which is why they pick up a freshly
# minted (but otherwise identical) text style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-apply the text format to the label/header cells so they get a new
# (content-identical) style entry, matching the re-export's style churn.
$ws.Range("A1").NumberFormat = "@"
$ws.Range("B1").NumberFormat = "@"
$ws.Range("C1").NumberFormat = "@"
$ws.Range("A2").NumberFormat = "@"

# The actual data change: B2's prediction score moves from 1 to the new
# computed value.
$ws.Range("B2").Value = -27.198268297674758
